# Elec Sources Used for Rlbty and Lst Cst Dsptch.xlsx
# Update "ESUfRaLCD-dispatch" sheet: insert "w ccs" rows (matching the
# ESUfRaLCD-reliability sheet), trim the unused placeholder rows/column,
# and move the active-sheet/selection state onto the dispatch tab.

$wb = $excel.ActiveWorkbook
$wsDispatch = $wb.Worksheets.Item("ESUfRaLCD-dispatch")
$wsReliability = $wb.Worksheets.Item("ESUfRaLCD-reliability")

# --- ESUfRaLCD-dispatch: insert 4 new rows (11-14) for the "w ccs" sources,
# pushing "small modular reactor" / hydrogen rows down to 15-17 ---
$wsDispatch.Rows.Item(11).Resize(4).Insert()

$wsDispatch.Range("A11").Value = "hard coal w ccs"
$wsDispatch.Range("B11").Value = "hard coal w ccs es"
$wsDispatch.Range("C11").Formula = '=IF(A11="","",CONCATENATE(A11," dispatch"))'

$wsDispatch.Range("A12").Value = "natural gas combined cycle w ccs"
$wsDispatch.Range("B12").Value = "natural gas combined cycle w ccs es"
$wsDispatch.Range("C12").Formula = '=IF(A12="","",CONCATENATE(A12," dispatch"))'

$wsDispatch.Range("A13").Value = "biomass w ccs"
$wsDispatch.Range("B13").Value = "biomass w CCS es"
$wsDispatch.Range("C13").Formula = '=IF(A13="","",CONCATENATE(A13," dispatch"))'

$wsDispatch.Range("A14").Value = "lignite w ccs"
$wsDispatch.Range("B14").Value = "lignite w CCS es"
$wsDispatch.Range("C14").Formula = '=IF(A14="","",CONCATENATE(A14," dispatch"))'

# --- Drop the now-unused placeholder rows (the original 18:82 block, shifted
# down by the 4 freshly-inserted rows) so the used range shrinks back to row 17 ---
$wsDispatch.Range("A18:A86").EntireRow.Delete()

# --- Drop the empty, styled-only column D ---
$wsDispatch.Columns.Item(4).Delete()

# --- Reliability sheet: clear old selection/active state, select A8:C11 ---
$wsReliability.Range("A8:C11").Select()

# --- Dispatch sheet becomes the active/visible tab, cursor on A4 ---
$wsDispatch.Activate()
$wsDispatch.Range("A4").Select()
